$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for specific rows to reflect repulled
# data / mean calculation changes.
$ws.Range("F2").Value = -7
$ws.Range("F5").Value = -9
$ws.Range("F8").Value = -3
$ws.Range("F13").Value = -2
$ws.Range("F17").Value = -3
